# Applies the "Actualizar 02-05-2021 21-30-53" availability-check update:
#   - refreshes the timestamp of the last existing 14-row availability block
#     (rows 800-813, column D) from 44232.87517093252 to 44232.8751709375
#   - appends one more 14-row availability block (rows 814-827) with status
#     "Disponible" and timestamp 44232.89639873675, including the per-row
#     hyperlink on column B (mirroring the very first block, rows 2-15)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. refresh timestamp on the previous block (rows 800-813) ---------
for ($r = 800; $r -le 813; $r++) {
    $ws.Range("D$r").Value() = 44232.8751709375
}

# --- 2. data for the new block (rows 814-827) ---------------------------
# columns: Nombre (A, shared string id as used elsewhere), URL (B, text +
# hyperlink), Disponibilidad (C, always "Disponible"), Fecha (D, serial date)
$names = @("Odoo","Blackbox","PowerBI","Dropbox","Odoo","GEE","UtilidadesOdoo","Filtros Dashboard","MapStore","GeoServer","Tomcat","Shiny","Github","EZ Exporter")

# text shown in column B (matches the other 56 rows using these same URLs)
$displayUrls  = @(
    "https://www.dataintelligence-group.com/",
    "https://serviciodashboard.azurewebsites.net/",
    "https://powerbi.microsoft.com/es-es/",
    "https://www.dropbox.com/",
    "https://dataintelligence.store/",
    "https://app-data-i.users.earthengine.app/",
    "https://odooutil.azurewebsites.net/",
    "https://filtradordashboard.azurewebsites.net/",
    "https://ide.dataintelligence-group.com/mapstore/#/",
    "https://ide.dataintelligence-group.com/geoserver/web/?0",
    "https://ide.dataintelligence-group.com/",
    "https://rpubs.com/dataintelligence/",
    "https://github.com/Sud-Austral/",
    "https://ezexporter.highviewapps.com/exports/export-profile/"
)
# actual hyperlink target (the MapStore entry splits off the "#/" fragment
# into SubAddress, exactly like every other occurrence of that link)
$linkAddresses  = @(
    "https://www.dataintelligence-group.com/",
    "https://serviciodashboard.azurewebsites.net/",
    "https://powerbi.microsoft.com/es-es/",
    "https://www.dropbox.com/",
    "https://dataintelligence.store/",
    "https://app-data-i.users.earthengine.app/",
    "https://odooutil.azurewebsites.net/",
    "https://filtradordashboard.azurewebsites.net/",
    "https://ide.dataintelligence-group.com/mapstore/",
    "https://ide.dataintelligence-group.com/geoserver/web/?0",
    "https://ide.dataintelligence-group.com/",
    "https://rpubs.com/dataintelligence/",
    "https://github.com/Sud-Austral/",
    "https://ezexporter.highviewapps.com/exports/export-profile/"
)
$subAddresses = @("","","","","","","","","/","","","","","")

$startRow = 814
$timestamp = 44232.89639873675

for ($i = 0; $i -lt $names.Count; $i++) {
    $r = $startRow + $i

    $ws.Range("A$r").Value() = $names[$i]
    $ws.Range("B$r").Value() = $displayUrls[$i]
    $ws.Range("C$r").Value() = "Disponible"
    $ws.Range("D$r").Value() = $timestamp

    if ($subAddresses[$i] -ne "") {
        $ws.Hyperlinks.Add($ws.Range("B$r"), $linkAddresses[$i], $subAddresses[$i]) | Out-Null
    } else {
        $ws.Hyperlinks.Add($ws.Range("B$r"), $linkAddresses[$i]) | Out-Null
    }
}
